# Auto-generated Excel COM-interop script
# Applies numeric cell updates (recalculated Teamcraft/Typhon profit figures)
# across the ALC/ARM/BSM/CUL/GSM/LTW/WVR sheets, matching the target diff.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 126
$ws.Range("I2").Value = 98.71429000000001
$ws.Range("J2").Value = 173.75
$ws.Range("K2").Value = 98.71429000000001
$ws.Range("L2").Value = 173.75
$ws.Range("M2").Value = 14.28570999999999
$ws.Range("N2").Value = -399.75
$ws.Range("H5").Value = 445.11765
$ws.Range("I5").Value = 82.625
$ws.Range("J5").Value = 767.3333
$ws.Range("K5").Value = 82.625
$ws.Range("L5").Value = 767.3333
$ws.Range("M5").Value = 32.375
$ws.Range("N5").Value = -997.3333
$ws.Range("H12").Value = 1416.8572
$ws.Range("I12").Value = 2700
$ws.Range("J12").Value = 903.6
$ws.Range("K12").Value = 2700
$ws.Range("L12").Value = 903.6
$ws.Range("M12").Value = -2530
$ws.Range("N12").Value = -1243.6
$ws.Range("H43").Value = 2250
$ws.Range("I43").Value = 1500
$ws.Range("K43").Value = 1500
$ws.Range("M43").Value = -1431
$ws.Range("H58").Value = 3210.4443
$ws.Range("I58").Value = 272
$ws.Range("K58").Value = 816
$ws.Range("M58").Value = -666
$ws.Range("H88").Value = 2202.2856
$ws.Range("J88").Value = 3004
$ws.Range("L88").Value = 3004
$ws.Range("N88").Value = -3816
$ws.Range("H91").Value = 2202.2856
$ws.Range("J91").Value = 3004
$ws.Range("L91").Value = 3004
$ws.Range("N91").Value = -5812
$ws.Range("H94").Value = 1005
$ws.Range("I94").Value = 1005
$ws.Range("K94").Value = 1005
$ws.Range("M94").Value = -554
$ws.Range("H112").Value = 1079.129
$ws.Range("I112").Value = 500
$ws.Range("J112").Value = 1098.4333
$ws.Range("K112").Value = 1500
$ws.Range("L112").Value = 3295.2999
$ws.Range("M112").Value = -392
$ws.Range("N112").Value = -5511.2999
$ws.Range("H137").Value = 86575
$ws.Range("I137").Value = 4466.5
$ws.Range("K137").Value = 13399.5
$ws.Range("M137").Value = -10849.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H32").Value = 16882.408
$ws.Range("I32").Value = 17395.922
$ws.Range("K32").Value = 17395.922
$ws.Range("M32").Value = -17108.922
$ws.Range("H45").Value = 4521.636
$ws.Range("I45").Value = 4646.2856
$ws.Range("J45").Value = 4303.5
$ws.Range("K45").Value = 4646.2856
$ws.Range("L45").Value = 4303.5
$ws.Range("M45").Value = -4269.2856
$ws.Range("N45").Value = -5057.5
$ws.Range("H61").Value = 3401.0527
$ws.Range("I61").Value = 3424.4
$ws.Range("K61").Value = 3424.4
$ws.Range("M61").Value = -3212.4
$ws.Range("H101").Value = 43333.332
$ws.Range("J101").Value = 43333.332
$ws.Range("L101").Value = 43333.332
$ws.Range("N101").Value = -49823.332
$ws.Range("H122").Value = 1936.409
$ws.Range("I122").Value = 1973
$ws.Range("K122").Value = 5919
$ws.Range("M122").Value = -3469
$ws.Range("H132").Value = 19023.834
$ws.Range("I132").Value = 2630.5625
$ws.Range("K132").Value = 7891.6875
$ws.Range("M132").Value = -5361.6875
$ws.Range("H136").Value = 3401.0527
$ws.Range("I136").Value = 3424.4
$ws.Range("K136").Value = 10273.2
$ws.Range("M136").Value = -7723.200000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3221.3572
$ws.Range("I20").Value = 3425
$ws.Range("J20").Value = 1999.5
$ws.Range("K20").Value = 3425
$ws.Range("L20").Value = 1999.5
$ws.Range("M20").Value = -3178
$ws.Range("N20").Value = -2493.5
$ws.Range("H134").Value = 122416
$ws.Range("I134").Value = 137530.5
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 412591.5
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = -410056.5
$ws.Range("N134").Value = -9570

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 178.75
$ws.Range("I13").Value = 178.75
$ws.Range("K13").Value = 536.25
$ws.Range("M13").Value = -368.25
$ws.Range("H34").Value = 934.5
$ws.Range("J34").Value = 925.1429000000001
$ws.Range("L34").Value = 2775.4287
$ws.Range("N34").Value = -2943.4287
$ws.Range("H55").Value = 2710.2632
$ws.Range("J55").Value = 2710.2632
$ws.Range("L55").Value = 8130.7896
$ws.Range("N55").Value = -8484.7896
$ws.Range("H131").Value = 734.89
$ws.Range("J131").Value = 734.89
$ws.Range("L131").Value = 2204.67
$ws.Range("N131").Value = -12284.67
$ws.Range("H139").Value = 2064.762
$ws.Range("I139").Value = 1454.909
$ws.Range("J139").Value = 2735.6
$ws.Range("K139").Value = 4364.727000000001
$ws.Range("L139").Value = 8206.799999999999
$ws.Range("M139").Value = 775.2729999999992
$ws.Range("N139").Value = -18486.8
$ws.Range("H140").Value = 1796.1765
$ws.Range("I140").Value = 1600
$ws.Range("K140").Value = 4800
$ws.Range("M140").Value = 380

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 4000
$ws.Range("J47").Value = 4000
$ws.Range("L47").Value = 4000
$ws.Range("N47").Value = -5136
$ws.Range("H70").Value = 10689.8
$ws.Range("I70").Value = 17571
$ws.Range("J70").Value = 4668.75
$ws.Range("K70").Value = 17571
$ws.Range("L70").Value = 4668.75
$ws.Range("M70").Value = -17301
$ws.Range("N70").Value = -5208.75
$ws.Range("H73").Value = 10689.8
$ws.Range("I73").Value = 17571
$ws.Range("J73").Value = 4668.75
$ws.Range("K73").Value = 17571
$ws.Range("L73").Value = 4668.75
$ws.Range("M73").Value = -16635
$ws.Range("N73").Value = -6540.75
$ws.Range("H122").Value = 1758.6
$ws.Range("I122").Value = 1509.5555
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 4528.666499999999
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -2078.666499999999
$ws.Range("N122").Value = -16900
$ws.Range("H132").Value = 99949
$ws.Range("I132").Value = 107543.6
$ws.Range("J132").Value = 87291.336
$ws.Range("K132").Value = 322630.8
$ws.Range("L132").Value = 261874.008
$ws.Range("M132").Value = -320100.8
$ws.Range("N132").Value = -266934.008

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3208.4167
$ws.Range("I22").Value = 2833.5
$ws.Range("K22").Value = 2833.5
$ws.Range("M22").Value = -2538.5
$ws.Range("H27").Value = 3208.4167
$ws.Range("I27").Value = 2833.5
$ws.Range("K27").Value = 2833.5
$ws.Range("M27").Value = -2726.5
$ws.Range("H40").Value = 3645.84
$ws.Range("I40").Value = 2238.2222
$ws.Range("J40").Value = 4437.625
$ws.Range("K40").Value = 2238.2222
$ws.Range("L40").Value = 4437.625
$ws.Range("M40").Value = -2102.2222
$ws.Range("N40").Value = -4709.625
$ws.Range("H122").Value = 936549.9
$ws.Range("I122").Value = 1155579
$ws.Range("K122").Value = 3466737
$ws.Range("M122").Value = -3464287
$ws.Range("H141").Value = 58500
$ws.Range("J141").Value = 58500
$ws.Range("L141").Value = 58500
$ws.Range("N141").Value = -68860

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 36233.332
$ws.Range("J109").Value = 36233.332
$ws.Range("L109").Value = 36233.332
$ws.Range("N109").Value = -39007.332
$ws.Range("H113").Value = 2080034.5
$ws.Range("I113").Value = 1302
$ws.Range("J113").Value = 9009143
$ws.Range("K113").Value = 3906
$ws.Range("L113").Value = 27027429
$ws.Range("M113").Value = -1736
$ws.Range("N113").Value = -27031769
$ws.Range("H122").Value = 1507
$ws.Range("I122").Value = 1582
$ws.Range("K122").Value = 4746
$ws.Range("M122").Value = -2296

